$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "last saved" date shown by the datetimeFigureOut
#    field on the Slide Master and on every Slide Layout (12 placeholders
#    total: 1 master + 11 layouts), from 21.10.2019 to 16.11.2020.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = "16.11.2020"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($l).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 2 ("Characteristics of Frontend / What to do in Frontend?")
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(1)
$text2 = $shape2.TextFrame.TextRange

# "...part for security checking)" -> split the trailing run into three runs:
# " security " / "checking done by " / "Backend)"
$hit = $text2.Find(" security checking)")
$start = $hit.Start
$hit.Text = " security checking done by Backend)"

$refreshed = $shape2.TextFrame.TextRange
$refreshed.Characters($start, 10).Text = " security "
$refreshed.Characters($start + 10, 17).Text = "checking done by "
$refreshed.Characters($start + 27, 8).Text = "Backend)"

# Extend the red warning line about request forging.
$text2b = $shape2.TextFrame.TextRange
$hit2 = $text2b.Find("Request forging might by-pass our Frontend checks totally!!! ")
$hit2.Text = "Request forging might by-pass our Frontend checks totally!!! (Or somebody writes new/changed Frontend without checks)"

# ---------------------------------------------------------------------------
# 3) Slide 4 ("Characteristics of Backend / What to do in Backend?")
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(1)
$text4 = $shape4.TextFrame.TextRange

# Collapse the three runs describing the single-backend bottleneck back into
# one run (same wording, now authored as a single run).
$hit3 = $text4.Find("Though often/mostly we only have one instance of the backend (then would be performance bottleneck)")
$hit3.Text = "Though often/mostly we only have one instance of the backend (then would be performance bottleneck)"

# "Protect Database" -> "Help Database"
$text4b = $shape4.TextFrame.TextRange
$hit4 = $text4b.Find("Protect Database by checking ")
$hit4.Text = "Help Database by checking "
